$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - Lefty O'Doul
$ws.Range("D9").Value = "2 Batting Titles"
$ws.Range("E9").Value = "11 year career"
$ws.Range("F9").Value = "Born on March 4, 1897"

# Row 10 - Ed Delahanty
$ws.Range("D10").Value = "2 Batting Titles"
$ws.Range("E10").Value = "16 year career"
$ws.Range("F10").Value = "Played for PHI, WSH, and CLE"

# Update the active cell selection to D11 (matches the diff)
$ws.Range("D11").Select()
